$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (2026-02-08 -> 2026-02-09, serial 46061 -> 46062) for every data row (2..361).
$ws.Range("C2:C361").Value = 46062
